# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" everywhere it
#   appears (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" stamps
#   tick forward to the new handoff run.
# - Columns E/F on Overview and column C on the language sheets are widened
#   to fit the new "Ready for handoff" status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-15 10:56:04"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-15 10:55:57"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-15 10:56:04"

# Widen the status columns so the longer "Ready for handoff" label fits.
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
